$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Resolving-Mac / Ccl12 / Ccr4 / ECs
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.23247666666667
$ws.Range("H2").Value = 60.69743
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04581866666666667
$ws.Range("N2").Value = 0.137456
$ws.Range("O2").Value = 0.4376255591461182
$ws.Range("P2").Value = 0.4376255591461182
$ws.Range("Q2").Value = 0.9270251042311112
$ws.Range("R2").Value = 8.34322593808
$ws.Range("S2").Value = 0.4376255591461182
$ws.Range("T2").Value = 0.4376255591461182

# Row 3: Resolving-Mac / Ccl12 / Ccr4 / FAPs
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 20.23247666666667
$ws.Range("H3").Value = 60.69743
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05887966666666666
$ws.Range("N3").Value = 0.176639
$ws.Range("O3").Value = 0.5623744408538818
$ws.Range("P3").Value = 0.5623744408538818
$ws.Range("Q3").Value = 1.191281481974444
$ws.Range("R3").Value = 10.72153333777
$ws.Range("S3").Value = 0.5623744408538818
$ws.Range("T3").Value = 0.5623744408538818
